$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.656.77'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '1.589.31'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("E6").Value = '  -2.27%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.62'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '1.811.40'
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("D13").Value = '1.588.69'
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("E15").Value = '  -4.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '26.631.15'
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("E18").Value = '  -2.51%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("E22").Value = '  -3.39%  '
$ws.Range("E23").Value = '  -3.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("E28").Value = '  -3.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0507'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.664'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +20.69%  '
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("D35").Value = '1.307.76'
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("E37").Value = '  -5.17%  '
$ws.Range("E38").Value = '  -3.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.831'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.792'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.34%  '
$ws.Range("D45").Value = '1.724.77'
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.836'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("E49").Value = '  -1.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0504'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.20%  '
